$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C fitness values change in three bands (spreadsheet rows, 1-based):
#   rows 2-16   (Generation 0-14)   : 7569 -> 7345
#   rows 17-82  (Generation 15-80)  : 7569 -> 7295
#   rows 83-252 (Generation 81-250) : 7569 -> 7293

$ws.Range("C2:C16").Value = 7345
$ws.Range("C17:C82").Value = 7295
$ws.Range("C83:C252").Value = 7293
